$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Absent" column (H) is the complement of the "Real" column (E):
# when a student was marked present for real (E=1) they are not absent (H=0),
# otherwise (E=0) they are absent (H=1). Fill in / correct column H for every
# data row so the consolidated report reflects this.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 3) { $lastRow = 21 }

for ($r = 3; $r -le $lastRow; $r++) {
    $realCell = $ws.Cells.Item($r, 5)   # column E = Real
    $absentCell = $ws.Cells.Item($r, 8) # column H = Absent

    $realValue = $realCell.Value2
    if ($realValue -eq 1) {
        $absentCell.Value = 0
    } else {
        $absentCell.Value = 1
    }
}
